# Adds 12 monthly "Reserva de Emergência / Reserva / Savings" entries
# (rows 80-91) to the LANCAMENTOS sheet, and updates the visible selection
# to G83 scrolled near the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LANCAMENTOS")

$dates = @(45818, 45848, 45879, 45910, 45940, 45971, 46001, 46032, 46063, 46091, 46122, 46152)

$startRow = 80
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i

    # Copy the formatting of the last existing data row (79) down onto the
    # new row so styles (date / currency number formats) match exactly.
    $ws.Range("A79:I79").Copy()
    $ws.Range(("A{0}:I{0}" -f $row)).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = "Reserva de Emergência"
    $ws.Cells.Item($row, 3).Value = "Reserva"
    $ws.Cells.Item($row, 4).Value = 500
    $ws.Cells.Item($row, 5).Value = "Savings"
    $ws.Cells.Item($row, 6).Value = "n"
    $ws.Cells.Item($row, 7).Value = "Unica"
    $ws.Cells.Item($row, 8).Value = 1
    $ws.Cells.Item($row, 9).Value = 1
}

$excel.CutCopyMode = $false

$ws.Activate()
$ws.Range("G83").Select()
